# Generate Report for Handback
# Adds "Latest Target File" (F) / "Latest Handback File" (G) data for the
# two localized sheets (zh-cn, de-de), marks the two source files as handed
# back (status text), and stamps the handback datetime per-language.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$statusHandedBack = "Handed back: in sync with en-US"

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# (This text is a shared string also surfaced on the Overview rollup sheet,
# so it is updated everywhere it is displayed.)
$wsOverview.Range("B2").Value = $statusHandedBack
$wsOverview.Range("C2").Value = $statusHandedBack
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

$wsZhCn.Range("C2").Value = $statusHandedBack
$wsZhCn.Range("C3").Value = $statusHandedBack

$wsDeDe.Range("C2").Value = $statusHandedBack
$wsDeDe.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet: populate F (Latest Target File) / G (Latest Handback File) ---
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/c24c424b7c4e9ae277b24f2f1857d8d583518275/e2e/adeda3a9-9890-4de8-905a-500ead680ab5.md", "", "", "adeda3a9-9890-4de8-905a-500ead680ab5.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d747c16a9c696bf072cc914f98a1c106fe55dfed/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/adeda3a9-9890-4de8-905a-500ead680ab5.665048fe12261c13e745c4e44c9a98b09469e35e.zh-cn.xlf", "", "", "adeda3a9-9890-4de8-905a-500ead680ab5.665048fe12261c13e745c4e44c9a98b09469e35e.zh-cn.xlf") | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/c24c424b7c4e9ae277b24f2f1857d8d583518275/e2e/fe5756af-a198-49bf-80ce-74ebceea6cc5.md", "", "", "fe5756af-a198-49bf-80ce-74ebceea6cc5.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d747c16a9c696bf072cc914f98a1c106fe55dfed/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/fe5756af-a198-49bf-80ce-74ebceea6cc5.64a45db4b3e3489f099f53c5e21ed4045082d4cc.zh-cn.xlf", "", "", "fe5756af-a198-49bf-80ce-74ebceea6cc5.64a45db4b3e3489f099f53c5e21ed4045082d4cc.zh-cn.xlf") | Out-Null

# zh-cn: Latest Handback DateTime (H) now populated (was the zero-date placeholder)
$wsZhCn.Range("H2").Value = "2016-03-18 02:55:45"
$wsZhCn.Range("H3").Value = "2016-03-18 02:55:45"

# --- de-de sheet: populate F (Latest Target File) / G (Latest Handback File) ---
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/c24c424b7c4e9ae277b24f2f1857d8d583518275/e2e/adeda3a9-9890-4de8-905a-500ead680ab5.md", "", "", "adeda3a9-9890-4de8-905a-500ead680ab5.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff82b5425bdda9bf8937b86ebda4a3e42ca90cf7/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/adeda3a9-9890-4de8-905a-500ead680ab5.665048fe12261c13e745c4e44c9a98b09469e35e.de-de.xlf", "", "", "adeda3a9-9890-4de8-905a-500ead680ab5.665048fe12261c13e745c4e44c9a98b09469e35e.de-de.xlf") | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/c24c424b7c4e9ae277b24f2f1857d8d583518275/e2e/fe5756af-a198-49bf-80ce-74ebceea6cc5.md", "", "", "fe5756af-a198-49bf-80ce-74ebceea6cc5.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ff82b5425bdda9bf8937b86ebda4a3e42ca90cf7/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/fe5756af-a198-49bf-80ce-74ebceea6cc5.64a45db4b3e3489f099f53c5e21ed4045082d4cc.de-de.xlf", "", "", "fe5756af-a198-49bf-80ce-74ebceea6cc5.64a45db4b3e3489f099f53c5e21ed4045082d4cc.de-de.xlf") | Out-Null

# de-de: Latest Handback DateTime (H) now populated (differs from zh-cn's stamp)
$wsDeDe.Range("H2").Value = "2016-03-18 02:55:59"
$wsDeDe.Range("H3").Value = "2016-03-18 02:55:59"
